# Qatar Stars League workbook update
# Author commit message: "Atualização de bases das ligas, do dia: 22-05-2024 às 20:16"
#
# The underlying source data feed re-sorted / re-paired a number of adjacent
# match rows (columns B..AB, i.e. everything except the running index in
# column A). For each of the row pairs below, the full record (id, teams,
# score, odds, P/L columns, etc.) held in row N was swapped with the record
# held in row N+1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs (first, second) whose B:AB contents must be exchanged.
$rowPairs = @(
    @(15, 16),
    @(24, 25),
    @(45, 46),
    @(60, 61),
    @(62, 63),
    @(68, 69),
    @(83, 84),
    @(90, 91),
    @(96, 97),
    @(104, 105),
    @(128, 129)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AB$rowA")
    $rangeB = $ws.Range("B$rowB`:AB$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
